$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37 update
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.24%  "

# Row 39 update
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'5.02"
$ws.Range("E39").Value = "  +1.99%  "

# Row 2
$ws.Range("D2").Value = "51.777.69"
$ws.Range("E2").Value = "  +5.80%  "

# Row 3
$ws.Range("D3").Value = "2.760.95"
$ws.Range("E3").Value = "  +4.65%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'117.99"
$ws.Range("E5").Value = "  +6.82%  "

# Row 6
$ws.Range("D6").Value = "'332.44"
$ws.Range("E6").Value = "  +3.24%  "

# Row 7
$ws.Range("D7").Value = "'0.534"
$ws.Range("E7").Value = "  +2.91%  "

# Row 8
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").Value = "'0.577"
$ws.Range("E9").Value = "  +7.01%  "

# Row 10
$ws.Range("D10").Value = "'41.84"
$ws.Range("E10").Value = "  +5.95%  "

# Row 11
$ws.Range("D11").Value = "'20.14"
$ws.Range("E11").Value = "  +1.44%  "

# Row 12
$ws.Range("D12").Value = "'0.0830"
$ws.Range("E12").Value = "  +2.68%  "

# Row 13
$ws.Range("E13").Value = "  +3.11%  "

# Row 14
$ws.Range("D14").Value = "'7.64"
$ws.Range("E14").Value = "  +6.04%  "

# Row 15
$ws.Range("D15").Value = "3.191.19"
$ws.Range("E15").Value = "  +4.67%  "

# Row 16
$ws.Range("D16").Value = "2.760.52"
$ws.Range("E16").Value = "  +4.61%  "

# Row 17
$ws.Range("D17").Value = "'0.885"
$ws.Range("E17").Value = "  +3.23%  "

# Row 18
$ws.Range("D18").Value = "51.698.26"
$ws.Range("E18").Value = "  +5.49%  "

# Row 19
$ws.Range("D19").Value = "'13.69"
$ws.Range("E19").Value = "  +6.72%  "

# Row 20
$ws.Range("D20").Value = "'2.99"
$ws.Range("E20").Value = "  +3.48%  "

# Row 21
$ws.Range("D21").Value = "'6.86"
$ws.Range("E21").Value = "  +3.00%  "

# Row 22
$ws.Range("D22").Value = "0.₃0964"
$ws.Range("E22").Value = "  +2.43%  "

# Row 23
$ws.Range("D23").Value = "'278.46"
$ws.Range("E23").Value = "  +3.23%  "

# Row 24
$ws.Range("D24").Value = "'69.69"
$ws.Range("E24").Value = "  -0.47%  "

# Row 25
$ws.Range("D25").Value = "'2.65"
$ws.Range("E25").Value = "  +3.92%  "

# Row 26
$ws.Range("D26").Value = "'26.89"
$ws.Range("E26").Value = "  +2.58%  "

# Row 27
$ws.Range("D27").Value = "'4.15"
$ws.Range("E27").Value = "  +0.52%  "

# Row 28
$ws.Range("E28").Value = "  +0.13%  "

# Row 29
$ws.Range("D29").Value = "'10.27"
$ws.Range("E29").Value = "  +2.09%  "

# Row 30
$ws.Range("D30").Value = "'2.23"
$ws.Range("E30").Value = "  +0.11%  "

# Row 31
$ws.Range("E31").Value = "  +2.21%  "

# Row 32
$ws.Range("D32").Value = "'35.60"
$ws.Range("E32").Value = "  +1.39%  "

# Row 33
$ws.Range("D33").Value = "'50.49"
$ws.Range("E33").Value = "  +2.38%  "

# Row 34
$ws.Range("D34").Value = "'5.61"
$ws.Range("E34").Value = "  +3.40%  "

# Row 35
$ws.Range("D35").Value = "'0.0824"
$ws.Range("E35").Value = "  +3.80%  "

# Row 36
$ws.Range("D36").Value = "'19.19"
$ws.Range("E36").Value = "  -0.01%  "

# Row 38
$ws.Range("E38").Value = "  +4.24%  "

# Row 40
$ws.Range("D40").Value = "'3.24"
$ws.Range("E40").Value = "  +3.09%  "

# Row 41
$ws.Range("D41").Value = "'130.50"
$ws.Range("E41").Value = "  +4.21%  "

# Row 42
$ws.Range("D42").Value = "'23.24"
$ws.Range("E42").Value = "  +2.85%  "

# Row 43
$ws.Range("E43").Value = "  +10.43%  "

# Row 45
$ws.Range("E45").Value = "  +4.56%  "

# Row 46
$ws.Range("D46").Value = "'2.37"
$ws.Range("E46").Value = "  +12.19%  "

# Row 47
$ws.Range("D47").Value = "2.117.82"
$ws.Range("E47").Value = "  +1.67%  "

# Row 48
$ws.Range("D48").Value = "'3.35"
$ws.Range("E48").Value = "  +4.36%  "

# Row 49
$ws.Range("E49").Value = "  +3.39%  "

# Row 50
$ws.Range("D50").Value = "'5.61"
$ws.Range("E50").Value = "  +8.27%  "

# Row 51
$ws.Range("D51").Value = "'9.01"
$ws.Range("E51").Value = "  +1.18%  "
